# Swap the two theme color schemes that live in this deck:
#   ppt/theme/theme2.xml (the theme actually applied to the slide master /
#   all slides) currently holds the "Integral" palette and needs to become
#   the default Office palette; ppt/theme/theme1.xml (used only by the
#   notes master) currently holds the Office palette and needs to become
#   Integral. The COM surface only round-trips color edits through the
#   slide-facing theme (ppt/theme/theme2.xml), so we rewrite its 12
#   scheme colors in place to match the "Office Theme" palette.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $tcs.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeColor 1  "000000"
Set-ThemeColor 2  "FFFFFF"
Set-ThemeColor 3  "44546A"
Set-ThemeColor 4  "E7E6E6"
Set-ThemeColor 5  "5B9BD5"
Set-ThemeColor 6  "ED7D31"
Set-ThemeColor 7  "A5A5A5"
Set-ThemeColor 8  "FFC000"
Set-ThemeColor 9  "4472C4"
Set-ThemeColor 10 "70AD47"
Set-ThemeColor 11 "0563C1"
Set-ThemeColor 12 "954F72"
